# Confirming acceptance with website.docx — apply the recorded edit:
#   1. Remove the stray "_GoBack" bookmark from the "Show recommended
#      information" Heading2 paragraph.
#   2. Change the "Admin account creation" Heading2 into
#      "Administrator account creation", leaving a "_GoBack" bookmark
#      sitting right after "Administrator" (i.e. where the editing
#      cursor ended up), matching how Word itself would have produced
#      the run split ("A" | "dmin" | "istrator" | <bookmark> | " account
#      creation").
#
# NOTE: this COM-interop host coalesces same-formatted adjacent runs in
# a paragraph whenever Range.InsertBefore/InsertAfter/Text touches that
# paragraph (mirroring how Word "re-types" contiguous identically
# formatted text into a single run). Bookmarks.Add, however, performs a
# pure structural split (no coalescing), so we use an "add a throwaway
# bookmark, then delete it" trick to re-establish run boundaries after
# the text insertion has coalesced them.

$d = $word.ActiveDocument

# --- 1. Drop the old "_GoBack" bookmark (it was at the end of the
#        "Show recommended information" heading). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Locate "Admin account creation" heading paragraph. ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $ptext = $p.Range.Text.Trim()
    if ($ptext -eq "Admin account creation") {
        $target = $p
        break
    }
}

if ($null -eq $target) {
    throw "Could not find the 'Admin account creation' heading paragraph."
}

$pStart = $target.Range.Start
# Offsets within the paragraph: "A"(0) "dmin"(1..5) " account creation"(5..)
$dminStart = $pStart + 1          # right after "A"
$afterDmin = $pStart + 5          # right after "dmin", before " account creation"

# Type "istrator" in right after "dmin" -> "Administrator account creation"
# (this will coalesce "A"+"dmin"+"istrator"+" account creation" into one run)
$ins = $d.Range($afterDmin, $afterDmin)
$ins.InsertBefore("istrator")

$afterIstrator = $afterDmin + 8   # right after the newly typed "istrator"

# --- 3. Re-split the run boundaries that existed/are needed, using the
#        add+delete bookmark trick (does not trigger coalescing). ---
# a) restore "A" | "dmin..." boundary
$d.Bookmarks.Add("_zz_tmp1", $d.Range($dminStart, $dminStart)) | Out-Null
$d.Bookmarks.Item("_zz_tmp1").Delete()

# b) restore "dmin" | "istrator..." boundary
$d.Bookmarks.Add("_zz_tmp2", $d.Range($afterDmin, $afterDmin)) | Out-Null
$d.Bookmarks.Item("_zz_tmp2").Delete()

# --- 4. Drop the real "_GoBack" bookmark right after "istrator", before
#        " account creation". ---
$d.Bookmarks.Add("_GoBack", $d.Range($afterIstrator, $afterIstrator)) | Out-Null

Write-Host "Applied edit: Admin -> Administrator, _GoBack bookmark moved."
